# Auto-generated script to apply "Add data for 2024-10-19" update
# Updates column K (2024 running totals) across Citywide Totals, By Neighborhood,
# and individual neighborhood detail sheets to reflect newly-recorded crime data.

$wb = $excel.ActiveWorkbook

# Citywide Totals
$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("K2").Value = 6495
$ws.Range("K3").Value = 6687
$ws.Range("K4").Value = 1392
$ws.Range("K5").Value = 481
$ws.Range("K6").Value = 7365
$ws.Range("K7").Value = 22420

# By Neighborhood
$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("K7").Value = 672
$ws.Range("K8").Value = 1474
$ws.Range("K9").Value = 101
$ws.Range("K11").Value = 413
$ws.Range("K19").Value = 657
$ws.Range("K20").Value = 536
$ws.Range("K22").Value = 69
$ws.Range("K23").Value = 224
$ws.Range("K26").Value = 32
$ws.Range("K29").Value = 1213
$ws.Range("K31").Value = 248
$ws.Range("K32").Value = 24
$ws.Range("K33").Value = 977
$ws.Range("K34").Value = 128
$ws.Range("K36").Value = 285
$ws.Range("K37").Value = 760
$ws.Range("K42").Value = 829
$ws.Range("K43").Value = 184
$ws.Range("K48").Value = 282
$ws.Range("K49").Value = 121
$ws.Range("K50").Value = 107
$ws.Range("K51").Value = 284
$ws.Range("K52").Value = 593
$ws.Range("K53").Value = 285
$ws.Range("K54").Value = 440
$ws.Range("K57").Value = 83
$ws.Range("K60").Value = 132
$ws.Range("K61").Value = 19
$ws.Range("K63").Value = 59
$ws.Range("K65").Value = 524
$ws.Range("K67").Value = 880
$ws.Range("K72").Value = 116
$ws.Range("K73").Value = 202
$ws.Range("K74").Value = 25
$ws.Range("K76").Value = 305
$ws.Range("K78").Value = 252
$ws.Range("K82").Value = 25
$ws.Range("K85").Value = 1037
$ws.Range("K86").Value = 137
$ws.Range("K88").Value = 238
$ws.Range("K89").Value = 333
$ws.Range("K90").Value = 210
$ws.Range("K93").Value = 83
$ws.Range("K96").Value = 240
$ws.Range("K97").Value = 179
$ws.Range("K98").Value = 113
$ws.Range("K99").Value = 373
$ws.Range("K101").Value = 22420

# West Ridge
$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("K6").Value = 102
$ws.Range("K7").Value = 240

# Auburn Gresham
$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("K2").Value = 221
$ws.Range("K3").Value = 220
$ws.Range("K6").Value = 181
$ws.Range("K7").Value = 672

# Belmont Cragin
$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("K6").Value = 138
$ws.Range("K7").Value = 413

# Uptown
$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("K5").Value = 2
$ws.Range("K6").Value = 97
$ws.Range("K7").Value = 333

# South Shore
$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("K2").Value = 338
$ws.Range("K7").Value = 1037

# Little Village
$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("K2").Value = 155
$ws.Range("K3").Value = 168
$ws.Range("K6").Value = 217
$ws.Range("K7").Value = 593

# Logan Square
$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("K6").Value = 121
$ws.Range("K7").Value = 285

# Austin
$ws = $wb.Worksheets.Item("Austin")
$ws.Range("K2").Value = 405
$ws.Range("K3").Value = 449
$ws.Range("K5").Value = 45
$ws.Range("K6").Value = 492
$ws.Range("K7").Value = 1474

# Garfield Park
$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("K3").Value = 349
$ws.Range("K7").Value = 977

# Grand Crossing
$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("K6").Value = 227
$ws.Range("K7").Value = 760

# New City
$ws = $wb.Worksheets.Item("New City")
$ws.Range("K3").Value = 128
$ws.Range("K7").Value = 524

# Woodlawn
$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("K2").Value = 97
$ws.Range("K3").Value = 154
$ws.Range("K7").Value = 373

# Gage Park
$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("K2").Value = 83
$ws.Range("K7").Value = 248

# North Lawndale
$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("K3").Value = 320
$ws.Range("K6").Value = 250
$ws.Range("K7").Value = 880

# Lincoln Park
$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("K2").Value = 26
$ws.Range("K7").Value = 121

# Loop
$ws = $wb.Worksheets.Item("Loop")
$ws.Range("K3").Value = 106
$ws.Range("K6").Value = 237
$ws.Range("K7").Value = 440

# Englewood
$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("K2").Value = 348
$ws.Range("K6").Value = 350
$ws.Range("K7").Value = 1213

# Lake View
$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("K2").Value = 44
$ws.Range("K7").Value = 282

# Chatham
$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("K2").Value = 194
$ws.Range("K3").Value = 198
$ws.Range("K6").Value = 215
$ws.Range("K7").Value = 657

# River North
$ws = $wb.Worksheets.Item("River North")
$ws.Range("K3").Value = 58
$ws.Range("K7").Value = 305

# Humboldt Park
$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("K3").Value = 253
$ws.Range("K7").Value = 829

# Rogers Park
$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("K2").Value = 76
$ws.Range("K6").Value = 87
$ws.Range("K7").Value = 252

# Douglas
$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("K3").Value = 78
$ws.Range("K4").Value = 14
$ws.Range("K7").Value = 224

# Chicago Lawn
$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("K2").Value = 187
$ws.Range("K7").Value = 536

# Grand Boulevard
$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("K3").Value = 86
$ws.Range("K7").Value = 285

# West Lawn
$ws = $wb.Worksheets.Item("West Lawn")
$ws.Range("K6").Value = 34
$ws.Range("K7").Value = 83

# Garfield Ridge
$ws = $wb.Worksheets.Item("Garfield Ridge")
$ws.Range("K3").Value = 34
$ws.Range("K7").Value = 128

# Wicker Park
$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Range("K6").Value = 66
$ws.Range("K7").Value = 113

# Lincoln Square
$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Range("K2").Value = 29
$ws.Range("K7").Value = 107

# East Village
$ws = $wb.Worksheets.Item("East Village")
$ws.Range("K4").Value = 1
$ws.Range("K7").Value = 32

# Avalon Park
$ws = $wb.Worksheets.Item("Avalon Park")
$ws.Range("K4").Value = 7
$ws.Range("K7").Value = 101

# Portage Park
$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("K2").Value = 68
$ws.Range("K6").Value = 69
$ws.Range("K7").Value = 202

# West Town
$ws = $wb.Worksheets.Item("West Town")
$ws.Range("K2").Value = 38
$ws.Range("K3").Value = 38
$ws.Range("K7").Value = 179

# United Center
$ws = $wb.Worksheets.Item("United Center")
$ws.Range("K3").Value = 73
$ws.Range("K7").Value = 238

# Galewood
$ws = $wb.Worksheets.Item("Galewood")
$ws.Range("K3").Value = 6
$ws.Range("K7").Value = 24

# Streeterville
$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("K2").Value = 23
$ws.Range("K7").Value = 137

# Washington Heights
$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("K6").Value = 52
$ws.Range("K7").Value = 210

# Little Italy, UIC
$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("K3").Value = 75
$ws.Range("K7").Value = 284

# Mckinley Park
$ws = $wb.Worksheets.Item("Mckinley Park")
$ws.Range("K3").Value = 18
$ws.Range("K7").Value = 83

# Morgan Park
$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Range("K4").Value = 11
$ws.Range("K7").Value = 132

# Hyde Park
$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("K4").Value = 26
$ws.Range("K7").Value = 184

# Clearing
$ws = $wb.Worksheets.Item("Clearing")
$ws.Range("K6").Value = 14
$ws.Range("K7").Value = 69

# Old Town
$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("K6").Value = 54
$ws.Range("K7").Value = 116

# Sheffield & DePaul
$ws = $wb.Worksheets.Item("Sheffield & DePaul")
$ws.Range("K5").Value = 13
$ws.Range("K6").Value = 25

# Mount Greenwood
$ws = $wb.Worksheets.Item("Mount Greenwood")
$ws.Range("K4").Value = 3
$ws.Range("K7").Value = 19

# Printers Row
$ws = $wb.Worksheets.Item("Printers Row")
$ws.Range("K2").Value = 4
$ws.Range("K7").Value = 25
